$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "LhfLJ648"
$ws.Range("B2").Value = 23092133
$ws.Range("C2").Value = "wrqraym19"
$ws.Range("D2").Value = "hQ!2`$kJ3"
$ws.Range("F2").Value = "TIktJnBp"
$ws.Range("G2").Value = "hiZP"
